# Update NATMI LR-pair output with new TPM-based values (Col11a1-Ddr1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.005705333333333333
$ws.Range("H2").Value = 0.017116
$ws.Range("I2").Value = 0.002541956540206457
$ws.Range("J2").Value = 0.002541956540206458
$ws.Range("M2").Value = 0.141694
$ws.Range("N2").Value = 0.425082
$ws.Range("O2").Value = 0.01763793963212447
$ws.Range("P2").Value = 0.01763793963212447
$ws.Range("Q2").Value = 0.0008084115013333332
$ws.Range("R2").Value = 0.007275703511999999
$ws.Range("S2").Value = 0.00004483487600364547
$ws.Range("T2").Value = 0.00004483487600364547
$ws.Range("G3").Value = 0.005705333333333333
$ws.Range("H3").Value = 0.017116
$ws.Range("I3").Value = 0.002541956540206457
$ws.Range("J3").Value = 0.002541956540206458
$ws.Range("O3").Value = 0.2714637835982539
$ws.Range("P3").Value = 0.2714637835982538
$ws.Range("Q3").Value = 0.01244218142444444
$ws.Range("R3").Value = 0.11197963282
$ws.Range("S3").Value = 0.000690049140146772
$ws.Range("T3").Value = 0.000690049140146772
$ws.Range("G4").Value = 0.005705333333333333
$ws.Range("H4").Value = 0.017116
$ws.Range("I4").Value = 0.002541956540206457
$ws.Range("J4").Value = 0.002541956540206458
$ws.Range("M4").Value = 5.710985666666667
$ws.Range("N4").Value = 17.132957
$ws.Range("O4").Value = 0.7108982767696218
$ws.Range("P4").Value = 0.7108982767696217
$ws.Range("Q4").Value = 0.03258307689022222
$ws.Range("R4").Value = 0.293247692012
$ws.Range("S4").Value = 0.00180707252405604
$ws.Range("T4").Value = 0.00180707252405604
$ws.Range("I5").Value = 0.3829116837922498
$ws.Range("J5").Value = 0.3829116837922499
$ws.Range("M5").Value = 0.141694
$ws.Range("N5").Value = 0.425082
$ws.Range("O5").Value = 0.01763793963212447
$ws.Range("P5").Value = 0.01763793963212447
$ws.Range("Q5").Value = 0.121776357808
$ws.Range("R5").Value = 1.095987220272
$ws.Range("S5").Value = 0.006753773163162836
$ws.Range("T5").Value = 0.006753773163162836
$ws.Range("I6").Value = 0.3829116837922498
$ws.Range("J6").Value = 0.3829116837922499
$ws.Range("O6").Value = 0.2714637835982539
$ws.Range("P6").Value = 0.2714637835982538
$ws.Range("S6").Value = 0.1039466544662223
$ws.Range("T6").Value = 0.1039466544662223
$ws.Range("I7").Value = 0.3829116837922498
$ws.Range("J7").Value = 0.3829116837922499
$ws.Range("M7").Value = 5.710985666666667
$ws.Range("N7").Value = 17.132957
$ws.Range("O7").Value = 0.7108982767696218
$ws.Range("P7").Value = 0.7108982767696217
$ws.Range("Q7").Value = 4.908203833474666
$ws.Range("R7").Value = 44.173834501272
$ws.Range("S7").Value = 0.2722112561628647
$ws.Range("T7").Value = 0.2722112561628647
$ws.Range("G8").Value = 1.379328
$ws.Range("H8").Value = 4.137983999999999
$ws.Range("I8").Value = 0.6145463596675437
$ws.Range("J8").Value = 0.6145463596675437
$ws.Range("M8").Value = 0.141694
$ws.Range("N8").Value = 0.425082
$ws.Range("O8").Value = 0.01763793963212447
$ws.Range("P8").Value = 0.01763793963212447
$ws.Range("Q8").Value = 0.195442501632
$ws.Range("R8").Value = 1.758982514688
$ws.Range("S8").Value = 0.01083933159295799
$ws.Range("T8").Value = 0.01083933159295799
$ws.Range("G9").Value = 1.379328
$ws.Range("H9").Value = 4.137983999999999
$ws.Range("I9").Value = 0.6145463596675437
$ws.Range("J9").Value = 0.6145463596675437
$ws.Range("O9").Value = 0.2714637835982539
$ws.Range("P9").Value = 0.2714637835982538
$ws.Range("Q9").Value = 3.00803620352
$ws.Range("R9").Value = 27.07232583168
$ws.Range("S9").Value = 0.1668270799918848
$ws.Range("T9").Value = 0.1668270799918848
$ws.Range("G10").Value = 1.379328
$ws.Range("H10").Value = 4.137983999999999
$ws.Range("I10").Value = 0.6145463596675437
$ws.Range("J10").Value = 0.6145463596675437
$ws.Range("M10").Value = 5.710985666666667
$ws.Range("N10").Value = 17.132957
$ws.Range("O10").Value = 0.7108982767696218
$ws.Range("P10").Value = 0.7108982767696217
$ws.Range("Q10").Value = 7.877322437631999
$ws.Range("R10").Value = 70.895901938688
$ws.Range("S10").Value = 0.436879948082701
$ws.Range("T10").Value = 0.4368799480827009
